$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D183:D235 values from "Yes" to "No"
$ws.Range("D183:D235").Value = "No"

# Update the view: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 175
$ws.Range("C181").Select()
